# This script applies a set of stock/quantity corrections to the
# CryCompanywiseStockReport workbook. For each affected item row, the
# Quantity (column F) and Value (column G = Rate * Qty, column D * F)
# are corrected, and the "Sub Total:" (column B) rows for each company
# group - as well as the overall Sub Total / Grand Total rows at the
# bottom of the report - are updated to reflect the corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F49").Value = 110
$ws.Range("G49").Value = 28680.3
$ws.Range("F53").Value = 75
$ws.Range("G53").Value = 22877.25
$ws.Range("B54").Value = 108130.53
$ws.Range("F78").Value = 7
$ws.Range("G78").Value = 6418.58
$ws.Range("B81").Value = 14914.14
$ws.Range("F106").Value = 41
$ws.Range("G106").Value = 10088.87
$ws.Range("B116").Value = 163748.26
$ws.Range("B160").Value = 57756
$ws.Range("F160").Value = 181
$ws.Range("G160").Value = 12025.64
$ws.Range("B161").Value = 53925
$ws.Range("F161").Value = 1
$ws.Range("G161").Value = 66.44
$ws.Range("F194").Value = 31
$ws.Range("G194").Value = 786.16
$ws.Range("B199").Value = 8095.88
$ws.Range("F208").Value = 80
$ws.Range("G208").Value = 5239.2
$ws.Range("F210").Value = 46
$ws.Range("G210").Value = 1853.8
$ws.Range("F211").Value = 203
$ws.Range("G211").Value = 13154.4
$ws.Range("F214").Value = 51
$ws.Range("G214").Value = 4502.79
$ws.Range("B221").Value = 52103.74
$ws.Range("F223").Value = 203
$ws.Range("G223").Value = 23202.9
$ws.Range("F224").Value = 2418
$ws.Range("G224").Value = 44733
$ws.Range("F225").Value = 40
$ws.Range("G225").Value = 2756.4
$ws.Range("B229").Value = 74691.10000000001
$ws.Range("F253").Value = 18
$ws.Range("G253").Value = 1476.18
$ws.Range("B266").Value = 98172.19
$ws.Range("F281").Value = 17
$ws.Range("G281").Value = 3263.32
$ws.Range("F284").Value = 85
$ws.Range("G284").Value = 12195.8
$ws.Range("B325").Value = 178710.23
$ws.Range("F362").Value = 53
$ws.Range("G362").Value = 8169.95
$ws.Range("F363").Value = 431
$ws.Range("G363").Value = 60594.29
$ws.Range("B365").Value = 76522.72
$ws.Range("F393").Value = 4
$ws.Range("G393").Value = 148.84
$ws.Range("B397").Value = 32809.7
$ws.Range("F404").Value = 92
$ws.Range("G404").Value = 16666.72
$ws.Range("F405").Value = 35
$ws.Range("G405").Value = 1527.75
$ws.Range("F410").Value = 125
$ws.Range("G410").Value = 4677.5
$ws.Range("F411").Value = 172
$ws.Range("G411").Value = 8763.4
$ws.Range("F412").Value = 11
$ws.Range("G412").Value = 2050.84
$ws.Range("F415").Value = 62
$ws.Range("G415").Value = 2061.5
$ws.Range("B421").Value = 111628.13
$ws.Range("F430").Value = 364
$ws.Range("G430").Value = 4786.6
$ws.Range("F431").Value = 475
$ws.Range("G431").Value = 6084.75
$ws.Range("F439").Value = 14
$ws.Range("G439").Value = 272.44
$ws.Range("F446").Value = 448
$ws.Range("G446").Value = 6599.04
$ws.Range("B447").Value = 78225.28
$ws.Range("F449").Value = 143
$ws.Range("G449").Value = 7245.81
$ws.Range("F452").Value = 377
$ws.Range("G452").Value = 7080.06
$ws.Range("F455").Value = 328
$ws.Range("G455").Value = 16426.24
$ws.Range("F458").Value = 226
$ws.Range("G458").Value = 10459.28
$ws.Range("F461").Value = 141
$ws.Range("G461").Value = 6825.81
$ws.Range("F463").Value = 343
$ws.Range("G463").Value = 3306.52
$ws.Range("F464").Value = 28
$ws.Range("G464").Value = 974.6799999999999
$ws.Range("B469").Value = 138344.87
$ws.Range("F509").Value = 47
$ws.Range("G509").Value = 2909.3
$ws.Range("B519").Value = 18634.6
$ws.Range("F522").Value = 109
$ws.Range("G522").Value = 2280.28
$ws.Range("F531").Value = 107
$ws.Range("G531").Value = 6806.27
$ws.Range("F534").Value = 31
$ws.Range("G534").Value = 2370.26
$ws.Range("F537").Value = 13
$ws.Range("G537").Value = 2242.63
$ws.Range("B538").Value = 62249.28
$ws.Range("F558").Value = 96
$ws.Range("G558").Value = 4352.64
$ws.Range("B564").Value = 5071.7
$ws.Range("F566").Value = 88
$ws.Range("G566").Value = 14122.24
$ws.Range("F575").Value = 173
$ws.Range("G575").Value = 9494.24
$ws.Range("F576").Value = 92
$ws.Range("G576").Value = 2520.8
$ws.Range("F580").Value = 78
$ws.Range("G580").Value = 5736.9
$ws.Range("F581").Value = 126
$ws.Range("G581").Value = 7999.74
$ws.Range("F582").Value = 131
$ws.Range("G582").Value = 9147.73
$ws.Range("F583").Value = 28
$ws.Range("G583").Value = 3969
$ws.Range("F584").Value = 11
$ws.Range("G584").Value = 1351.13
$ws.Range("F586").Value = 158
$ws.Range("G586").Value = 13957.72
$ws.Range("B588").Value = 142018.9
$ws.Range("F592").Value = 137
$ws.Range("G592").Value = 17885.35
$ws.Range("B596").Value = 38863.44
$ws.Range("F608").Value = 115
$ws.Range("G608").Value = 4965.7
$ws.Range("B611").Value = 12403.03
$ws.Range("F614").Value = 2
$ws.Range("G614").Value = 164.08
$ws.Range("B621").Value = 5721.17
$ws.Range("F669").Value = 42
$ws.Range("G669").Value = 3425.52
$ws.Range("F672").Value = 52
$ws.Range("G672").Value = 6786
$ws.Range("F675").Value = 303
$ws.Range("G675").Value = 11295.84
$ws.Range("F676").Value = 31
$ws.Range("G676").Value = 2811.08
$ws.Range("F677").Value = 170
$ws.Range("G677").Value = 11825.2
$ws.Range("F680").Value = 81
$ws.Range("G680").Value = 10935.81
$ws.Range("F681").Value = 201
$ws.Range("G681").Value = 24262.71
$ws.Range("F682").Value = 10
$ws.Range("G682").Value = 1207.1
$ws.Range("B683").Value = 85233
$ws.Range("F692").Value = 154
$ws.Range("G692").Value = 23163.14
$ws.Range("F693").Value = 5
$ws.Range("G693").Value = 139.15
$ws.Range("F705").Value = 57
$ws.Range("G705").Value = 1884.42
$ws.Range("F706").Value = 32
$ws.Range("G706").Value = 561.28
$ws.Range("F707").Value = 64
$ws.Range("G707").Value = 2751.36
$ws.Range("B713").Value = 73940.06
$ws.Range("F734").Value = 0
$ws.Range("G734").Value = 0
$ws.Range("B744").Value = 58637.69
$ws.Range("F751").Value = 2227
$ws.Range("G751").Value = 363245.97
$ws.Range("F752").Value = 254
$ws.Range("G752").Value = 71848.98
$ws.Range("F753").Value = 351
$ws.Range("G753").Value = 50772.15
$ws.Range("F754").Value = 58
$ws.Range("G754").Value = 2212.12
$ws.Range("F757").Value = 115
$ws.Range("G757").Value = 7762.5
$ws.Range("B759").Value = 523902.07
$ws.Range("B764").Value = 3219427.17
$ws.Range("B765").Value = 3219427.17
